$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New data row appended by the Adafruit IO sync (row 70)
$row = 70

$ws.Cells.Item($row, 1).Value = "2024-09-25T18:06:40Z"
$ws.Cells.Item($row, 2).Value = "temperature"

# "25" looks numeric, so force it to stay text like the rest of the
# Value column, then drop back to the default style so no new cell
# style gets attached to the cell.
$ws.Cells.Item($row, 3).NumberFormat = "@"
$ws.Cells.Item($row, 3).Value = "25"
$ws.Cells.Item($row, 3).Style = "Normal"

$ws.Cells.Item($row, 4).Value = "N/A"
$ws.Cells.Item($row, 5).Value = "N/A"
$ws.Cells.Item($row, 6).Value = "N/A"
